# "this added last report 18-12-24"
# Update the BL Audit Form with the latest report figures for 18-12-2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / report date labels -------------------------------------------------
$ws.Range("B1").Value = "18.12.2024"
$ws.Range("F34").Value = "19.12.2024 payment "

# --- Stock table (rows 9-21) ------------------------------------------------------
$ws.Range("C9").Value = 442423
$ws.Range("C10").Value = 1500
$ws.Range("C12").Value = 2000
$ws.Range("C13").Value = 42
$ws.Range("C16").Value = 70
$ws.Range("C17").Value = 78

# --- Hand cash / bank rows ---------------------------------------------------------
$ws.Range("E22").Value = 18833
$ws.Range("E23").Value = 11329

# --- Credit section ------------------------------------------------------------
$ws.Range("E27").Value = 22500

# --- Bank guarantee section ------------------------------------------------------
$ws.Range("E34").Value = 100000

# --- Scroll the sheet view down so row 30 is at the top (matches the saved view) --
# (selection/active cell stays on F35, only the window's scroll position moves)
$excel.ActiveWindow.ScrollRow = 30
